# new_ph_pl.xlsx — "break out stock.yaml completed"
#
# 1) "day" sheet: append 10 new data rows (90-99) — a fresh stock.yaml batch
#    dated 08/07/2024, pushing the sheet dimension from A1:I89 to A1:I99.
# 2) "week" sheet: normalize bsecode column D for rows 94-96 (Colgate,
#    United Breweries, National Aluminium) from text to numeric, matching
#    the numeric bsecode storage used everywhere else in the workbook.

$wb = $excel.ActiveWorkbook

# --- 1) Append new rows to the "day" sheet -------------------------------
$day = $wb.Worksheets.Item("day")

$newRows = @(
    @(1,  "MRF",        "Mrf Limited",                                     500290, -0.04, 128900.05, 20120,   "day", "08/07/2024 11:34:44"),
    @(2,  "ABBOTINDIA",  "Abbott India Limited",                            500488, 0.35,  27921.7,   6982,    "day", "08/07/2024 11:34:44"),
    @(3,  "BRITANNIA",   "Britannia Industries Limited",                    500825, 0.39,  5568.55,   171210,  "day", "08/07/2024 11:34:44"),
    @(4,  "LALPATHLAB",  "Dr. Lal Path Labs Ltd.",                          539524, 0.74,  2906.4,    362484,  "day", "08/07/2024 11:34:44"),
    @(5,  "MUTHOOTFIN",  "Muthoot Finance Limited",                         533398, -1.45, 1784.4,    163811,  "day", "08/07/2024 11:34:44"),
    @(6,  "GODREJCP",    "Godrej Consumer Products Limited",                532424, 3.8,   1426,      2471652, "day", "08/07/2024 11:34:44"),
    @(7,  "TATACONSUM",  "TATA Consumer Products Ltd",                      500800, 1.18,  1150.8,    1958676, "day", "08/07/2024 11:34:44"),
    @(8,  "AUBANK",      "AU Small Finance Bank",                           540611, -4.45, 642.7,     6191001, "day", "08/07/2024 11:34:44"),
    @(9,  "ABFRL",       "Aditya Birla Fashion And Retail Limited",         535755, -1.65, 322.25,    3706077, "day", "08/07/2024 11:34:44"),
    @(10, "M&MFIN",      "Mahindra & Mahindra Financial Services Limited",  532720, 0.77,  302.8,     1916926, "day", "08/07/2024 11:34:44")
)

$startRow = 90
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $day.Cells.Item($r, 1).Value = $data[0]
    $day.Cells.Item($r, 2).Value = $data[1]
    $day.Cells.Item($r, 3).Value = $data[2]
    $day.Cells.Item($r, 4).Value = $data[3]
    $day.Cells.Item($r, 5).Value = $data[4]
    $day.Cells.Item($r, 6).Value = $data[5]
    $day.Cells.Item($r, 7).Value = $data[6]
    $day.Cells.Item($r, 8).Value = $data[7]
    $day.Cells.Item($r, 9).Value = $data[8]
}

# --- 2) Fix bsecode type on the "week" sheet for rows 94-96 --------------
$week = $wb.Worksheets.Item("week")

$week.Cells.Item(94, 4).Value = 500830
$week.Cells.Item(95, 4).Value = 532478
$week.Cells.Item(96, 4).Value = 532234
